# Refresh the coin-price table with the latest snapshot values (cryptos list
# auto-update, GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds text like "23.898.60" / "1.000" in the source
# data (grouped-digits, not locale numbers). Pre-format the cells we are
# about to rewrite as Text so Excel does not silently coerce values such as
# "1.000" into the number 1.
$ws.Range("D2:D21").NumberFormat = "@"
$ws.Range("D23:D47").NumberFormat = "@"
$ws.Range("D49:D51").NumberFormat = "@"

$ws.Range("D2").Value = "23.898.60"
$ws.Range("E2").Value = "  -0.41%  "

$ws.Range("D3").Value = "1.648.98"
$ws.Range("E3").Value = "  -0.05%  "

$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.75%  "

$ws.Range("D5").Value = "310.62"
$ws.Range("E5").Value = "  +0.23%  "

$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.78%  "

$ws.Range("D7").Value = "0.3895"
$ws.Range("E7").Value = "  -0.78%  "

$ws.Range("D8").Value = "0.3838"

$ws.Range("D9").Value = "50.97"
$ws.Range("E9").Value = "  -0.39%  "

$ws.Range("D10").Value = "1.343"
$ws.Range("E10").Value = "  -1.38%  "

$ws.Range("D11").Value = "1.000"
$ws.Range("E11").Value = "  -0.72%  "

$ws.Range("D12").Value = "0.08444"
$ws.Range("E12").Value = "  -0.39%  "

$ws.Range("D13").Value = "23.83"
$ws.Range("E13").Value = "  -0.56%  "

$ws.Range("D14").Value = "7.023"
$ws.Range("E14").Value = "  -2.37%  "

$ws.Range("D15").Value = "7.917"

$ws.Range("D16").Value = "0.00001315"
$ws.Range("E16").Value = "  +0.18%  "

$ws.Range("D17").Value = "1.647.37"
$ws.Range("E17").Value = "  -0.31%  "

$ws.Range("D18").Value = "93.92"
$ws.Range("E18").Value = "  -0.65%  "

$ws.Range("D19").Value = "0.06966"
$ws.Range("E19").Value = "  -0.37%  "

$ws.Range("D20").Value = "19.52"
$ws.Range("E20").Value = "  -2.83%  "

$ws.Range("D21").Value = "6.938"
$ws.Range("E21").Value = "  +0.41%  "

$ws.Range("E22").Value = "  -0.75%  "

$ws.Range("D23").Value = "13.63"
$ws.Range("E23").Value = "  -0.39%  "

$ws.Range("D24").Value = "23.904.06"

$ws.Range("D25").Value = "2.445"
$ws.Range("E25").Value = "  -2.54%  "

$ws.Range("D26").Value = "2.907"
$ws.Range("E26").Value = "  -5.11%  "

$ws.Range("D27").Value = "21.97"
$ws.Range("E27").Value = "  -1.30%  "

$ws.Range("D28").Value = "154.19"
$ws.Range("E28").Value = "  -0.75%  "

$ws.Range("D29").Value = "5.391"
$ws.Range("E29").Value = "  +1.14%  "

$ws.Range("D30").Value = "137.22"
$ws.Range("E30").Value = "  -1.83%  "

$ws.Range("D31").Value = "7.720"
$ws.Range("E31").Value = "  -1.76%  "

$ws.Range("D32").Value = "2.482"
$ws.Range("E32").Value = "  -1.25%  "

$ws.Range("D33").Value = "1.830.30"
$ws.Range("E33").Value = "  -0.07%  "

$ws.Range("D34").Value = "0.08118"
$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("D35").Value = "0.9897"
$ws.Range("E35").Value = "  -4.55%  "

$ws.Range("D36").Value = "0.02921"
$ws.Range("E36").Value = "  -3.30%  "

$ws.Range("D37").Value = "6.683"
$ws.Range("E37").Value = "  -0.49%  "

$ws.Range("D38").Value = "0.2680"
$ws.Range("E38").Value = "  -1.38%  "

$ws.Range("D39").Value = "10.45"
$ws.Range("E39").Value = "  -4.20%  "

$ws.Range("D40").Value = "0.09116"
$ws.Range("E40").Value = "  -0.48%  "

$ws.Range("D41").Value = "0.7551"
$ws.Range("E41").Value = "  +0.07%  "

$ws.Range("D42").Value = "13.41"
$ws.Range("E42").Value = "  -0.89%  "

$ws.Range("D43").Value = "1.420"
$ws.Range("E43").Value = "  -0.70%  "

$ws.Range("D44").Value = "16.71"
$ws.Range("E44").Value = "  +2.74%  "

$ws.Range("D45").Value = "0.6930"
$ws.Range("E45").Value = "  -0.08%  "

$ws.Range("D46").Value = "2.437"
$ws.Range("E46").Value = "  -1.82%  "

$ws.Range("D47").Value = "4.097"
$ws.Range("E47").Value = "  +0.07%  "

$ws.Range("E48").Value = "  -0.76%  "

$ws.Range("D49").Value = "0.08267"
$ws.Range("E49").Value = "  -0.03%  "

$ws.Range("D50").Value = "134.36"
$ws.Range("E50").Value = "  -0.02%  "

$ws.Range("B51").Value = "Flow"
$ws.Range("C51").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D51").Value = "1.222"
$ws.Range("E51").Value = "  -0.62%  "
